# Weekly data refresh: a new week's worth of price observations (rows for
# variety "Magnum" and "Sin especificar") is inserted right after row 80,
# pushing the existing rows 81-190 down to 83-192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows at row 81 (existing rows 81-190 shift to 83-192).
$ws.Rows.Item(81).Resize(2).Insert()

# New row 81: Poroto verde / Magnum
$ws.Range("A81").Value = 2
$ws.Range("B81").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C81").Value = "Coquimbo"
$ws.Range("D81").Value = 44671
$ws.Range("E81").Value = 4
$ws.Range("F81").Value = 100112031
$ws.Range("G81").Value = "Poroto verde"
$ws.Range("H81").Value = "Magnum"
$ws.Range("I81").Value = "Primera"
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 13000
$ws.Range("L81").Value = 15000
$ws.Range("M81").Value = 14000
$ws.Range("N81").Value = "$/malla 25 kilos"
$ws.Range("O81").Value = "Provincia de Limarí"
$ws.Range("P81").Value = 560
$ws.Range("Q81").Value = 25
$ws.Range("R81").Value = "Hortaliza"

# New row 82: Poroto verde / Sin especificar
$ws.Range("A82").Value = 2
$ws.Range("B82").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C82").Value = "Coquimbo"
$ws.Range("D82").Value = 44671
$ws.Range("E82").Value = 4
$ws.Range("F82").Value = 100112031
$ws.Range("G82").Value = "Poroto verde"
$ws.Range("H82").Value = "Sin especificar"
$ws.Range("I82").Value = "Primera"
$ws.Range("J82").Value = 500
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 21000
$ws.Range("M82").Value = 20500
$ws.Range("N82").Value = "$/malla 25 kilos"
$ws.Range("O82").Value = "Provincia de Limarí"
$ws.Range("P82").Value = 820
$ws.Range("Q82").Value = 25
$ws.Range("R82").Value = "Hortaliza"
